$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (A6) down to the new
# row (A7) so the new cell picks up the same style index instead of
# Excel minting a brand-new cellXf.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 0.3293225559127213
